$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("K6").Value = 2

$ws.Range("K5").Select()
